# Katie Meso TA data up to 2/17/2020
# Append six new data rows (90-95) to Sheet1, mirroring the existing
# CRM-accuracy log, and move the visible selection to just past the new
# last row (F96), matching the author's save-time cursor position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- new row data --------------------------------------------------
# row, A(date serial), B(CRM value), C(batch value), E(batch #), F(note)
$rows = @(
    @{ r = 90; a = 43860; b = 2235.9172055679701; c = 2235.0700000000002; e = 155; f = "CRM opened 1/10/2020 (Dudgeon)" },
    @{ r = 91; a = 43864; b = 2248.65524442177;   c = 2235.0700000000002; e = 155; f = "CRM opened 1/10/2020 (Dudgeon)" },
    @{ r = 92; a = 43867; b = 2246.27558568522;   c = 2235.0700000000002; e = 155; f = "CRM opened 1/10/2020 (Dudgeon)" },
    @{ r = 93; a = 43871; b = 2245.1589170981301; c = 2235.0700000000002; e = 155; f = "CRM opened 1/10/2020 (Dudgeon)" },
    @{ r = 94; a = 43874; b = 2241.65662945434;   c = 2235.0700000000002; e = 155; f = "CRM opened 1/10/2020 (Dudgeon)" },
    @{ r = 95; a = 43878; b = 2217.49579116116;   c = 2207.0300000000002; e = 169; f = "CRM opened 2/17/2020 (Silbiger bottle for Dudgeon)" }
)

foreach ($row in $rows) {
    $r = $row.r

    $ws.Range("A$r").Value = $row.a
    # Copy the date's number format (style) down from the row above so the
    # new cells share the existing date style instead of minting a new one.
    $ws.Range("A" + ($r - 1)).Copy()
    $ws.Range("A$r").PasteSpecial(-4122) # xlPasteFormats

    $ws.Range("B$r").Value = $row.b
    $ws.Range("C$r").Value = $row.c
    $ws.Range("D$r").Formula = "=100*(B$r-C$r)/C$r"
    $ws.Range("E$r").Value = $row.e
    $ws.Range("F$r").Value = $row.f
}

$excel.CutCopyMode = $false

# ---- move selection to mirror the saved view state ------------------
$ws.Range("F96").Select()
